# Updated cryptos list on Wed Jul  3 08:56:27 UTC 2024 with GitHub Actions
#
# The "Price" column (D) and "Volume(1h)" column (E) are stored as plain
# text cells in the workbook (prices use '.' as both thousands- and
# decimal-separator, so they can't round-trip as numbers; percentages keep
# padding whitespace). When a price string looks like an ordinary decimal
# number (e.g. "27.63"), assigning it straight to Range.Value would make
# Excel auto-coerce the cell to a real number, which both changes the
# stored type and would silently drop a trailing zero (e.g. "22.80" ->
# 22.8). Set-TextValue forces such values to stay text (via a leading
# apostrophe) and then reapplies a plain/default cell style so no stray
# "quote prefix" / number-format style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cleanStyle = $ws.Range("B2").Style
    $cell.Value = "'" + $value
    $cell.Style = $cleanStyle
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.774.20"
$ws.Range("E2").Value = "  -3.05%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.346.90"
$ws.Range("E3").Value = "  -2.81%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "566.79"
$ws.Range("E5").Value = "  -2.22%  "

# Row 6 - Solana
$ws.Range("E6").Value = "  -1.02%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.40%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  -0.78%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.38%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +1.10%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.26%  "

# Row 14 - Avalanche
Set-TextValue "D14" "27.63"
$ws.Range("E14").Value = "  -2.11%  "

# Row 15 - WrappedEther
$ws.Range("E15").Value = "  -2.80%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -1.50%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "60.763.77"
$ws.Range("E17").Value = "  -3.12%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  -1.53%  "

# Row 19 - Chainlink
Set-TextValue "D19" "14.49"

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -2.00%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "375.22"
$ws.Range("E21").Value = "  -3.03%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -0.29%  "

# Row 23 - Litecoin
Set-TextValue "D23" "74.80"
$ws.Range("E23").Value = "  -0.63%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.11%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.496.60"
$ws.Range("E25").Value = "  -2.36%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -5.74%  "

# Row 27 - Kaspa
Set-TextValue "D27" "0.174"
$ws.Range("E27").Value = "  -4.11%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("E28").Value = "  +0.19%  "

# Row 29 - RenderToken
Set-TextValue "D29" "7.32"
$ws.Range("E29").Value = "  -4.05%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -1.26%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -3.16%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "22.80"

# Row 34 - Fetch.AI
$ws.Range("E34").Value = "  -3.96%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  -0.45%  "

# Row 36 - Monero
Set-TextValue "D36" "168.57"
$ws.Range("E36").Value = "  -0.53%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -4.95%  "

# Row 38 - Aptos
$ws.Range("E38").Value = "  -2.35%  "

# Row 40 - RenzoRestakedETH
$ws.Range("E40").Value = "  -2.62%  "

# Row 41 - Hedera
$ws.Range("E41").Value = "  -3.02%  "

# Row 42 - Mantle
$ws.Range("E42").Value = "  -3.53%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  -1.54%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  -5.46%  "

# Row 45 - ONDO
$ws.Range("E45").Value = "  -3.89%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.459.76"
$ws.Range("E46").Value = "  -4.19%  "

# Row 47 - Cosmos
$ws.Range("E47").Value = "  -4.05%  "

# Row 48 - FirstDigitalUSD
$ws.Range("E48").Value = "  +0.08%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "22.23"
$ws.Range("E49").Value = "  -1.63%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  -2.03%  "

# Row 51 - SuiNetwork
$ws.Range("E51").Value = "  +0.40%  "
